$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Het Patel"
$ws.Range("B3").Value = "hetpatel5542@gmail.com"
$ws.Range("C3").Value = "GCET"

# Force the phone number to be stored as text (matching the existing
# Phone Number column in row 2), instead of being auto-converted to a
# numeric value.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "7698545581"
$ws.Range("D3").Style = $ws.Range("D2").Style

$ws.Range("E3").Value = "STATIC_COMBO"
$ws.Range("F3").Value = "OFFLINE"
